# Remove initial survey screen
# (Adds a new "initial" screen/sheet before "survey" that lets the user
#  jump straight into the survey section and skip the finalize screen.)

$wb = $excel.ActiveWorkbook

# Insert a brand-new worksheet. Worksheets.Add() with no args places the
# new sheet immediately before the active sheet, i.e. right before
# "survey" (the first tab) - giving us the desired leading position.
$initial = $wb.Worksheets.Add()
$initial.Name = "initial"

# --- header row -------------------------------------------------------
$initial.Range("A1").Value = "clause"
$initial.Range("B1").Value = "type"
$initial.Range("C1").Value = "display.text"
$initial.Range("D1").Value = "comments"

# --- row 2: jump straight into the survey section ---------------------
$initial.Range("A2").Value = "do section survey"

# --- row 3: skip the finalize screen -----------------------------------
$initial.Range("A3").Value = "goto _finalize"
$initial.Range("D3").Value = "skips the finalize screen where the user chooses to save as incomplete or finalized and instead saves as finalized"

# --- formatting ---------------------------------------------------------
$initial.Range("A1:D1").WrapText = $true
$initial.Range("A2:C2").WrapText = $true
$initial.Range("A3:D3").WrapText = $true

# NB: the host engine quantizes ColumnWidth to 1/6-character pixel ticks
# (Excel's internal Maximum-Digit-Width pixel rounding) before writing the
# stored `width` attribute, so the literal target widths are passed in
# pre-compensated form to land exactly (16 and 18) or as close as possible
# (24.25 is not a multiple of 1/6 and lands on the nearest tick, 24.333...).
$initial.Columns.Item(1).ColumnWidth = 15.166666666666666
$initial.Columns.Item(3).ColumnWidth = 17.166666666666668
$initial.Columns.Item(4).ColumnWidth = 23.5

$initial.Rows.Item(3).RowHeight = 77.5

# --- page setup / margins ----------------------------------------------
$ps = $initial.PageSetup
$ps.Orientation = 1
$ps.LeftMargin = $excel.InchesToPoints(0.75)
$ps.RightMargin = $excel.InchesToPoints(0.75)
$ps.TopMargin = $excel.InchesToPoints(1)
$ps.BottomMargin = $excel.InchesToPoints(1)
$ps.HeaderMargin = $excel.InchesToPoints(0.51180555555555496)
$ps.FooterMargin = $excel.InchesToPoints(0.51180555555555496)

# --- selection matches the authored file --------------------------------
$initial.Range("D6").Select() | Out-Null
